# Tweaking colorado gas model
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Modbus register descriptions (rows 2-8, column C) now include p/t offsets
$ws.Range("C2").Value = "(5000, 0), p:(5100, 2), t:(5200, 4)"
$ws.Range("C3").Value = "(5001, 0), p:(5101, 2), t:(5201, 4)"
$ws.Range("C4").Value = "(5003, 0), p:(5103, 2), t:(5203, 4)"
$ws.Range("C5").Value = "(5004, 0), p:(5104, 2), t:(5204, 4)"
$ws.Range("C6").Value = "(5005, 0), p:(5105, 2), t:(5205, 4)"
$ws.Range("C7").Value = "(5006, 0), p:(5106, 2), t:(5206, 4)"
$ws.Range("C8").Value = "(5007, 0), p:(5107 , 2), t:(5207, 4)"

# Load point parameter specs (column C) simplified to just p/t (no pt/tt)
$ws.Range("C21").Value = "p: 8000, t: 9000"
$ws.Range("C22").Value = "p: 8001, t: 9001"
$ws.Range("C24").Value = "p: 8002, t: 9002"
$ws.Range("C26").Value = "p: 8003, t: 9003"
$ws.Range("C28").Value = "p: 8004, t: 9004"

# Update window/selection state to match
$ws.Range("C11").Select()
